$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the two new worksheets (as copies of "CreateProfile", which already
#    carries the exact column widths / styles / header row used by the new
#    sheets) and place them at the end of the workbook.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("CreateProfile")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "AddDeleteCCNotPrepopulated"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "AddDeleteCCPrepopulated"

# ---------------------------------------------------------------------------
# 2) AddDeleteCCNotPrepopulated (sheet4): wipe the copied rows 2:3 and the
#    extra-wide "ProfileName" column, then lay out the CC columns AI:AU.
# ---------------------------------------------------------------------------
$ws4.Range("A2:AH3").Clear()
$ws4.Columns.Item(3).ColumnWidth = 27.59

$ws4.Range("AI1").Value = "TitleCC"
$ws4.Range("AJ1").Value = "FirstNameCC"
$ws4.Range("AK1").Value = "MiddleNameCC"
$ws4.Range("AL1").Value = "LastNameCC"
$ws4.Range("AM1").Value = "SuffixCC"
$ws4.Range("AN1").Value = "NickNameCC"
$ws4.Range("AO1").Value = "CardNumberCC"
$ws4.Range("AP1").Value = "ExpMonthCC"
$ws4.Range("AQ1").Value = "ExpYearCC"
$ws4.Range("AR1").Value = "AL1CC"
$ws4.Range("AS1").Value = "AL2CC"
$ws4.Range("AT1").Value = "ZIPCC"
$ws4.Range("AU1").Value = "ZIPExtCC"

$ws4.Range("C2").Value = "All Fields CC"
$ws4.Range("D2").Value = "Y"
$ws4.Range("E2").Value = "742"
$ws4.Range("F2").Value = "a_Access AutoNoCFtp"
$ws4.Range("G2").Value = "4249"
$ws4.Range("H2").Value = "a_Access AutoNoCFtpDemo"
$ws4.Range("P2").Value = "15 Toledo Road"
$ws4.Range("R2").Value = "United States"
$ws4.Range("S2").Value = "22201"
$ws4.Range("AI2").Value = "Mr."
$ws4.Range("AJ2").Value = "Timothy"
$ws4.Range("AK2").Value = "Mac"
$ws4.Range("AL2").Value = "Siefert"
$ws4.Range("AM2").Value = "Sr."
$ws4.Range("AN2").Value = "Tim MasterCard"
$ws4.Range("AO2").Value = "5146312200000035"
$ws4.Range("AP2").Value = "12- DEC"
$ws4.Range("AQ2").Value = "2028"
$ws4.Range("AR2").Value = "365 Kanyakumari"
$ws4.Range("AS2").Value = "Room 8"
$ws4.Range("AT2").Value = "22201"
$ws4.Range("AU2").Value = "1234"

$ws4.Range("C3").Value = "Required Fields CC"
$ws4.Range("D3").Value = "Y"
$ws4.Range("E3").Value = "742"
$ws4.Range("F3").Value = "a_Access AutoNoCFtp"
$ws4.Range("G3").Value = "4249"
$ws4.Range("H3").Value = "a_Access AutoNoCFtpDemo"
$ws4.Range("P3").Value = "15 Toledo Road"
$ws4.Range("R3").Value = "United States"
$ws4.Range("S3").Value = "22201"
$ws4.Range("AJ3").Value = "Brandon"
$ws4.Range("AL3").Value = "McCulumm"
$ws4.Range("AN3").Value = "Brandon Amex"
$ws4.Range("AO3").Value = "371449635392376"
$ws4.Range("AP3").Value = "1 - JAN"
$ws4.Range("AQ3").Value = "2028"
$ws4.Range("AR3").Value = "2225 Kendsha road"
$ws4.Range("AT3").Value = "21093"

$ws4.Range("Y1").Select()
$ws4.Range("I3:AH3").Select()
$ws4.Range("I3").Activate()

# ---------------------------------------------------------------------------
# 3) AddDeleteCCPrepopulated (sheet5): wipe the copied rows 2:3, widen the
#    "ProfileName" column, and write the pre-populated CC profile.
# ---------------------------------------------------------------------------
$ws5.Range("A2:AH3").Clear()
$ws5.Columns.Item(3).ColumnWidth = 27.59

$ws5.Range("AI1").Value = "NickNameCC"
$ws5.Range("AJ1").Value = "CardNumberCC"
$ws5.Range("AK1").Value = "ExpMonthCC"
$ws5.Range("AL1").Value = "ExpYearCC"

$ws5.Range("C2").Value = "Required Fields CC"
$ws5.Range("D2").Value = "Y"
$ws5.Range("E2").Value = "742"
$ws5.Range("F2").Value = "a_Access AutoNoCFtp"
$ws5.Range("G2").Value = "4249"
$ws5.Range("H2").Value = "a_Access AutoNoCFtpDemo"
$ws5.Range("J2").Value = "Delta Corp"
$ws5.Range("K2").Value = "Mr."
$ws5.Range("L2").Value = "Ross"
$ws5.Range("M2").Value = "KT"
$ws5.Range("N2").Value = "Evan"
$ws5.Range("O2").Value = "Sr."
$ws5.Range("P2").Value = "256987 Nolm Ct"
$ws5.Range("Q2").Value = "Suite 678"
$ws5.Range("R2").Value = "United States"
$ws5.Range("S2").Value = "21054"
$ws5.Range("T2").Value = "1234"
$ws5.Range("W2").Value = "iahmed@govolution.com"
$ws5.Range("W2").Style = $template.Range("W3").Style
$ws5.Range("X2").Value = "240"
$ws5.Range("Y2").Value = "628"
$ws5.Range("Z2").Value = "0790"
$ws5.Range("AA2").Value = "240"
$ws5.Range("AB2").Value = "628"
$ws5.Range("AC2").Value = "0791"
$ws5.Range("AD2").Value = "410"
$ws5.Range("AE2").Value = "628"
$ws5.Range("AF2").Value = "0792"
$ws5.Range("AG2").Value = "123"
$ws5.Range("AH2").Value = "This Profile is for Add CC"
$ws5.Range("AI2").Value = "Ross Discover"
$ws5.Range("AJ2").Value = "6011000993026909"
$ws5.Range("AK2").Value = "1 - JAN"
$ws5.Range("AL2").Value = "2028"

$ws5.Range("AE1").Select()
$ws5.Range("AK2").Select()
$ws5.Range("AK2").Activate()

# ---------------------------------------------------------------------------
# 4) Workbook-level view bookkeeping to mirror the target bookViews entry.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayWorkbookTabs = $true
$wb.Windows.Item(1).ScrollWorkbookTabs(1, 1)
$ws5.Activate()
